# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several leve rows
# on the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1921
$ws.Range("I19").Value = 2359.6
$ws.Range("J19").Value = 824.5
$ws.Range("K19").Value = 2359.6
$ws.Range("L19").Value = 824.5
$ws.Range("M19").Value = -2184.6
$ws.Range("N19").Value = -1174.5
$ws.Range("H107").Value = 567.6429000000001
$ws.Range("J107").Value = 410
$ws.Range("L107").Value = 410
$ws.Range("N107").Value = -4250
$ws.Range("H112").Value = 1650.1034
$ws.Range("I112").Value = 1825
$ws.Range("J112").Value = 1583.4762
$ws.Range("K112").Value = 5475
$ws.Range("L112").Value = 4750.4286
$ws.Range("M112").Value = -4367
$ws.Range("N112").Value = -6966.4286
$ws.Range("H129").Value = 2149.1875
$ws.Range("I129").Value = 1959.1333
$ws.Range("K129").Value = 5877.3999
$ws.Range("M129").Value = -877.3999000000003
$ws.Range("H138").Value = 4936.7
$ws.Range("I138").Value = 3349.75
$ws.Range("J138").Value = 5113.028
$ws.Range("K138").Value = 10049.25
$ws.Range("L138").Value = 15339.084
$ws.Range("M138").Value = -4909.25
$ws.Range("N138").Value = -25619.084

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3669.875
$ws.Range("I74").Value = 1842.5714
$ws.Range("K74").Value = 1842.5714
$ws.Range("M74").Value = -968.5714
$ws.Range("H77").Value = 3669.875
$ws.Range("I77").Value = 1842.5714
$ws.Range("K77").Value = 9212.857
$ws.Range("M77").Value = -4844.857
$ws.Range("H132").Value = 2953.6924
$ws.Range("I132").Value = 2680.0625
$ws.Range("J132").Value = 4204.5713
$ws.Range("K132").Value = 8040.1875
$ws.Range("L132").Value = 12613.7139
$ws.Range("M132").Value = -5510.1875
$ws.Range("N132").Value = -17673.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 146.39131
$ws.Range("I7").Value = 105.85714
$ws.Range("J7").Value = 209.44444
$ws.Range("K7").Value = 105.85714
$ws.Range("L7").Value = 209.44444
$ws.Range("M7").Value = 7.142859999999999
$ws.Range("N7").Value = -435.44444
$ws.Range("H62").Value = 2085.7144
$ws.Range("I62").Value = 1721
$ws.Range("J62").Value = 2997.5
$ws.Range("K62").Value = 1721
$ws.Range("L62").Value = 2997.5
$ws.Range("M62").Value = -1097
$ws.Range("N62").Value = -4245.5
$ws.Range("H65").Value = 2085.7144
$ws.Range("I65").Value = 1721
$ws.Range("J65").Value = 2997.5
$ws.Range("K65").Value = 8605
$ws.Range("L65").Value = 14987.5
$ws.Range("M65").Value = -5485
$ws.Range("N65").Value = -21227.5
$ws.Range("H68").Value = 93739.57000000001
$ws.Range("J68").Value = 93739.57000000001
$ws.Range("L68").Value = 93739.57000000001
$ws.Range("N68").Value = -95237.57000000001
$ws.Range("H71").Value = 93739.57000000001
$ws.Range("J71").Value = 93739.57000000001
$ws.Range("L71").Value = 281218.71
$ws.Range("N71").Value = -288706.71
$ws.Range("H94").Value = 1388.7222
$ws.Range("I94").Value = 912.25
$ws.Range("J94").Value = 1524.8572
$ws.Range("K94").Value = 912.25
$ws.Range("L94").Value = 1524.8572
$ws.Range("M94").Value = -461.25
$ws.Range("N94").Value = -2426.8572
$ws.Range("H122").Value = 2347.4443
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 2546.7144
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 7640.1432
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -12540.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 68.3
$ws.Range("I2").Value = 49
$ws.Range("J2").Value = 97.25
$ws.Range("K2").Value = 294
$ws.Range("L2").Value = 583.5
$ws.Range("M2").Value = -181
$ws.Range("N2").Value = -809.5
$ws.Range("H117").Value = 3858.8333
$ws.Range("I117").Value = 1495.4546
$ws.Range("J117").Value = 7572.7144
$ws.Range("K117").Value = 4486.3638
$ws.Range("L117").Value = 22718.1432
$ws.Range("M117").Value = -1044.3638
$ws.Range("N117").Value = -29602.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 51428.637
$ws.Range("I122").Value = 67716.87
$ws.Range("J122").Value = 16525.285
$ws.Range("K122").Value = 203150.61
$ws.Range("L122").Value = 49575.855
$ws.Range("M122").Value = -200700.61
$ws.Range("N122").Value = -54475.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2886.375
$ws.Range("I7").Value = 2316.2856
$ws.Range("K7").Value = 2316.2856
$ws.Range("M7").Value = -2204.2856
$ws.Range("H61").Value = 5143.1113
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H68").Value = 3581.0588
$ws.Range("J68").Value = 4677.3335
$ws.Range("L68").Value = 4677.3335
$ws.Range("N68").Value = -6175.3335
$ws.Range("H71").Value = 3581.0588
$ws.Range("J71").Value = 4677.3335
$ws.Range("L71").Value = 23386.6675
$ws.Range("N71").Value = -30874.6675
$ws.Range("H100").Value = 7356.5713
$ws.Range("I100").Value = 2681.4707
$ws.Range("J100").Value = 27225.75
$ws.Range("K100").Value = 2681.4707
$ws.Range("L100").Value = 27225.75
$ws.Range("M100").Value = -2140.4707
$ws.Range("N100").Value = -28307.75
$ws.Range("H113").Value = 5143.1113
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 46505.22
$ws.Range("I122").Value = 2718.8823
$ws.Range("K122").Value = 8156.646900000001
$ws.Range("M122").Value = -5706.646900000001
$ws.Range("H126").Value = 2886.375
$ws.Range("I126").Value = 2316.2856
$ws.Range("K126").Value = 6948.8568
$ws.Range("M126").Value = -4478.8568
$ws.Range("H136").Value = 5002.5884
$ws.Range("I136").Value = 2503.6667
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 7511.000100000001
$ws.Range("L136").Value = 33000
$ws.Range("M136").Value = -4961.000100000001
$ws.Range("N136").Value = -38100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H122").Value = 1855.4445
$ws.Range("I122").Value = 1363.8182
$ws.Range("K122").Value = 4091.4546
$ws.Range("M122").Value = -1641.4546
$ws.Range("H136").Value = 7142650.5
$ws.Range("I136").Value = 8187367.5
$ws.Range("K136").Value = 24562102.5
$ws.Range("M136").Value = -24559552.5
